{"js": "// Append six new numbered list items (ListParagraph / numId 1) to the end\n// of the document body, right after the \"Test now passes.\" item, describing\n// the work done to require an email address when registering a user.\n\nconst body = context.document.body;\n\nconst newItems = [\n  \"Make user an instance variable in the controller so that information gets saved and user doesn\u2019t have to fill out all forms again. \",\n  [\n    \"View now expects @user instance variable, so need to make sure it\u2019s available in the /users/new get \",\n    \"route. Placed @user = User.new into get request now all is fine \",\n  ],\n  \"Wrote test for requiring an email address. Not passing as it says \u2018wrong number of arguments (1 for 0) \",\n  \"Made changes to sign up method so it takes user as an argument\",\n  \"Had to change all the fill in properties in the test. So username became user.username. Password became user.password. \",\n  \"Test passes. Had to include validates presence of email in user model and change the tests so the sign up method takes one parameter  (user). \",\n];\n\nfor (const item of newItems) {\n  const runs = Array.isArray(item) ? item : [item];\n\n  // Append a new paragraph at the end of the body, seeded with the first\n  // run's text. Inserting at the body's \"End\" (rather than \"After\" a given\n  // paragraph) makes the new paragraph pick up the preceding ListParagraph\n  // style and numbered-list (numId 1) formatting automatically.\n  const newPara = body.insertParagraph(runs[0], \"End\");\n  await context.sync();\n\n  // Any additional runs get appended to the end of this same paragraph.\n  for (let i = 1; i < runs.length; i++) {\n    const tail = newPara.getRange(\"End\");\n    tail.insertText(runs[i], \"End\");\n    await context.sync();\n  }\n}\n", "ps1": "# Append six new numbered list items (ListParagraph / numId 1) to the end\n# of the document body, right after the \"Test now passes.\" item, describing\n# the work done to require an email address when registering a user.\n\n$d = $word.ActiveDocument\n\n# Anchor on the current last paragraph in the body: \"Test now passes.\"\n$n = $d.Paragraphs.Count\n$anchor = $d.Paragraphs.Item($n)\n\n$newItems = @(\n    ,@(\"Make user an instance variable in the controller so that information gets saved and user doesn\u2019t have to fill out all forms again. \")\n    ,@(\"View now expects @user instance variable, so need to make sure it\u2019s available in the /users/new get \", \"route. Placed @user = User.new into get request now all is fine \")\n    ,@(\"Wrote test for requiring an email address. Not passing as it says \u2018wrong number of arguments (1 for 0) \")\n    ,@(\"Made changes to sign up method so it takes user as an argument\")\n    ,@(\"Had to change all the fill in properties in the test. So username became user.username. Password became user.password. \")\n    ,@(\"Test passes. Had to include validates presence of email in user model and change the tests so the sign up method takes one parameter  (user). \")\n)\n\nforeach ($runs in $newItems) {\n    # Insert a new paragraph after the anchor; it inherits the anchor's\n    # ListParagraph style + numbered-list (numId 1) formatting automatically.\n    $anchor.Range.InsertParagraphAfter()\n\n    $count = $d.Paragraphs.Count\n    $newPara = $d.Paragraphs.Item($count)\n    $newPara.Range.Text = $runs[0]\n\n    for ($i = 1; $i -lt $runs.Length; $i++) {\n        $tail = $newPara.Range\n        $tail.Collapse(0)\n        $tail.InsertAfter($runs[$i])\n    }\n\n    $anchor = $newPara\n}\n\nWrite-Output \"done\"\n"}
